$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 135 (ALC)
$ws.Range("H135").Value = 2414.1
$ws.Range("I135").Value = 2486.8948
$ws.Range("K135").Value = 22382.0532
$ws.Range("M135").Value = -19847.0532

# Row 137 (ALC)
$ws.Range("H137").Value = 1383.3928
$ws.Range("I137").Value = 975.1111
$ws.Range("K137").Value = 2925.3333
$ws.Range("M137").Value = -375.3332999999998

# Row 141 (ALC)
$ws.Range("H141").Value = 3439.75
$ws.Range("I141").Value = 3298
$ws.Range("K141").Value = 9894
$ws.Range("M141").Value = -4714

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (ARM)
$ws.Range("H61").Value = 4242.6665
$ws.Range("I61").Value = 2187.375
$ws.Range("K61").Value = 2187.375
$ws.Range("M61").Value = -1975.375

# Row 63 (ARM)
$ws.Range("H63").Value = 4661.5
$ws.Range("I63").Value = 6186
$ws.Range("J63").Value = 2374.75
$ws.Range("K63").Value = 6186
$ws.Range("L63").Value = 2374.75
$ws.Range("M63").Value = -5500
$ws.Range("N63").Value = -3746.75

# Row 66 (ARM)
$ws.Range("H66").Value = 4661.5
$ws.Range("I66").Value = 6186
$ws.Range("J66").Value = 2374.75
$ws.Range("K66").Value = 30930
$ws.Range("L66").Value = 11873.75
$ws.Range("M66").Value = -27498
$ws.Range("N66").Value = -18737.75

# Row 74 (ARM)
$ws.Range("H74").Value = 2218.4614
$ws.Range("I74").Value = 2323.1904
$ws.Range("J74").Value = 1778.6
$ws.Range("K74").Value = 2323.1904
$ws.Range("L74").Value = 1778.6
$ws.Range("M74").Value = -1449.1904
$ws.Range("N74").Value = -3526.6

# Row 77 (ARM)
$ws.Range("H77").Value = 2218.4614
$ws.Range("I77").Value = 2323.1904
$ws.Range("J77").Value = 1778.6
$ws.Range("K77").Value = 11615.952
$ws.Range("L77").Value = 8893
$ws.Range("M77").Value = -7247.951999999999
$ws.Range("N77").Value = -17629

# Row 86 (ARM)
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Row 89 (ARM)
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Row 136 (ARM)
$ws.Range("H136").Value = 4242.6665
$ws.Range("I136").Value = 2187.375
$ws.Range("K136").Value = 6562.125
$ws.Range("M136").Value = -4012.125

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 3234.4736
$ws.Range("I31").Value = 1273.4
$ws.Range("J31").Value = 5413.4443
$ws.Range("K31").Value = 1273.4
$ws.Range("L31").Value = 5413.4443
$ws.Range("M31").Value = -978.4000000000001
$ws.Range("N31").Value = -6003.4443

# Row 34 (CRP)
$ws.Range("H34").Value = 3234.4736
$ws.Range("I34").Value = 1273.4
$ws.Range("J34").Value = 5413.4443
$ws.Range("K34").Value = 1273.4
$ws.Range("L34").Value = 5413.4443
$ws.Range("M34").Value = -1071.4
$ws.Range("N34").Value = -5817.4443

# Row 132 (CRP)
$ws.Range("H132").Value = 2076.5715
$ws.Range("I132").Value = 2102
$ws.Range("K132").Value = 6306
$ws.Range("M132").Value = -3776

# Row 134 (CRP)
$ws.Range("H134").Value = 3364.62
$ws.Range("I134").Value = 2121.5588
$ws.Range("K134").Value = 6364.676399999999
$ws.Range("M134").Value = -3829.676399999999

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 9027.875
$ws.Range("I5").Value = 4582.5386
$ws.Range("J5").Value = 14281.454
$ws.Range("K5").Value = 13747.6158
$ws.Range("L5").Value = 42844.362
$ws.Range("M5").Value = -13635.6158
$ws.Range("N5").Value = -43068.362

# Row 131 (CUL)
$ws.Range("H131").Value = 5632
$ws.Range("I131").Value = 3318.625
$ws.Range("J131").Value = 7945.375
$ws.Range("K131").Value = 9955.875
$ws.Range("L131").Value = 23836.125
$ws.Range("M131").Value = -4915.875
$ws.Range("N131").Value = -33916.125

# Row 135 (CUL)
$ws.Range("H135").Value = 9027.875
$ws.Range("I135").Value = 4582.5386
$ws.Range("J135").Value = 14281.454
$ws.Range("K135").Value = 41242.8474
$ws.Range("L135").Value = 128533.086
$ws.Range("M135").Value = -38707.8474
$ws.Range("N135").Value = -133603.086

# Row 137 (CUL)
$ws.Range("H137").Value = 1740.1428
$ws.Range("I137").Value = 1845
$ws.Range("J137").Value = 1111
$ws.Range("K137").Value = 5535
$ws.Range("L137").Value = 3333
$ws.Range("M137").Value = -435
$ws.Range("N137").Value = -13533

$ws = $wb.Worksheets.Item("GSM")
# Row 126 (GSM)
$ws.Range("H126").Value = 8000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 24000
$ws.Range("N126").Value = -28940
$ws.Range("M126").ClearContents()

# Row 132 (GSM)
$ws.Range("H132").Value = 2931.0322
$ws.Range("I132").Value = 1699.9546
$ws.Range("J132").Value = 5940.3335
$ws.Range("K132").Value = 5099.8638
$ws.Range("L132").Value = 17821.0005
$ws.Range("M132").Value = -2569.8638
$ws.Range("N132").Value = -22881.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (LTW)
$ws.Range("H22").Value = 839.5
$ws.Range("I22").Value = 796.36365
$ws.Range("K22").Value = 796.36365
$ws.Range("M22").Value = -501.36365

# Row 27 (LTW)
$ws.Range("H27").Value = 839.5
$ws.Range("I27").Value = 796.36365
$ws.Range("K27").Value = 796.36365
$ws.Range("M27").Value = -689.36365

# Row 36 (LTW)
$ws.Range("H36").Value = 58611.5
$ws.Range("J36").Value = 58611.5
$ws.Range("L36").Value = 58611.5
$ws.Range("N36").Value = -59735.5

# Row 74 (LTW)
$ws.Range("H74").Value = 35000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

# Row 77 (LTW)
$ws.Range("H77").Value = 35000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

# Row 81 (LTW)
$ws.Range("H81").Value = 52500
$ws.Range("J81").Value = 52500
$ws.Range("L81").Value = 52500
$ws.Range("N81").Value = -54496

# Row 82 (LTW)
$ws.Range("H82").Value = 2581.6667
$ws.Range("I82").Value = 2461.7778
$ws.Range("J82").Value = 2761.5
$ws.Range("K82").Value = 2461.7778
$ws.Range("L82").Value = 2761.5
$ws.Range("M82").Value = -2100.7778
$ws.Range("N82").Value = -3483.5

# Row 84 (LTW)
$ws.Range("H84").Value = 52500
$ws.Range("J84").Value = 52500
$ws.Range("L84").Value = 157500
$ws.Range("N84").Value = -167484

# Row 85 (LTW)
$ws.Range("H85").Value = 2581.6667
$ws.Range("I85").Value = 2461.7778
$ws.Range("J85").Value = 2761.5
$ws.Range("K85").Value = 2461.7778
$ws.Range("L85").Value = 2761.5
$ws.Range("M85").Value = -1213.7778
$ws.Range("N85").Value = -5257.5

# Row 124 (LTW)
$ws.Range("H124").Value = 186665.67
$ws.Range("J124").Value = 186665.67
$ws.Range("L124").Value = 186665.67
$ws.Range("N124").Value = -196485.67

# Row 132 (LTW)
$ws.Range("H132").Value = 4185.2
$ws.Range("I132").Value = 3376.5881
$ws.Range("K132").Value = 10129.7643
$ws.Range("M132").Value = -7599.764299999999

# Row 136 (LTW)
$ws.Range("H136").Value = 5332.7827
$ws.Range("I136").Value = 3386
$ws.Range("K136").Value = 10158
$ws.Range("M136").Value = -7608
